{"js": "// Update the worksheet date and every \"a\u00d7b=c\" answer cell to the new\n// values. Every <w:t> run in the document is a distinct, unique string,\n// so a simple ordered find/replace (old -> new) fully reproduces the\n// diff while leaving all run/paragraph formatting untouched.\nconst replacements = [\n  [\"2024-06-22 Saturday\", \"2024-06-23 Sunday\"],\n  [\"11\u00d790=990\", \"16\u00d737=592\"],\n  [\"67\u00d780=5360\", \"85\u00d759=5015\"],\n  [\"78\u00d785=6630\", \"52\u00d732=1664\"],\n  [\"51\u00d743=2193\", \"11\u00d711=121\"],\n  [\"21\u00d758=1218\", \"52\u00d792=4784\"],\n  [\"93\u00d746=4278\", \"24\u00d765=1560\"],\n  [\"64\u00d757=3648\", \"35\u00d726=910\"],\n  [\"92\u00d797=8924\", \"95\u00d736=3420\"],\n  [\"14\u00d715=210\", \"99\u00d794=9306\"],\n  [\"64\u00d768=4352\", \"22\u00d773=1606\"],\n  [\"63\u00d752=3276\", \"30\u00d724=720\"],\n  [\"33\u00d743=1419\", \"86\u00d768=5848\"],\n  [\"40\u00d725=1000\", \"33\u00d764=2112\"],\n  [\"85\u00d719=1615\", \"76\u00d741=3116\"],\n  [\"63\u00d711=693\", \"45\u00d774=3330\"],\n  [\"30\u00d790=2700\", \"88\u00d739=3432\"],\n  [\"49\u00d733=1617\", \"56\u00d746=2576\"],\n  [\"34\u00d761=2074\", \"84\u00d782=6888\"],\n  [\"92\u00d765=5980\", \"42\u00d794=3948\"],\n  [\"78\u00d761=4758\", \"62\u00d793=5766\"],\n  [\"37\u00d741=1517\", \"46\u00d712=552\"],\n  [\"56\u00d786=4816\", \"65\u00d772=4680\"],\n  [\"38\u00d713=494\", \"71\u00d750=3550\"],\n  [\"45\u00d766=2970\", \"39\u00d773=2847\"],\n  [\"93\u00d723=2139\", \"81\u00d748=3888\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every \"a\u00d7b=c\" answer cell to the new\n# values. Every text run in the document is a distinct, unique string,\n# so a simple ordered Find/Replace (old -> new) across the whole\n# document body fully reproduces the diff while leaving all\n# run/paragraph formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-06-22 Saturday\", \"2024-06-23 Sunday\"),\n  @(\"11\u00d790=990\", \"16\u00d737=592\"),\n  @(\"67\u00d780=5360\", \"85\u00d759=5015\"),\n  @(\"78\u00d785=6630\", \"52\u00d732=1664\"),\n  @(\"51\u00d743=2193\", \"11\u00d711=121\"),\n  @(\"21\u00d758=1218\", \"52\u00d792=4784\"),\n  @(\"93\u00d746=4278\", \"24\u00d765=1560\"),\n  @(\"64\u00d757=3648\", \"35\u00d726=910\"),\n  @(\"92\u00d797=8924\", \"95\u00d736=3420\"),\n  @(\"14\u00d715=210\", \"99\u00d794=9306\"),\n  @(\"64\u00d768=4352\", \"22\u00d773=1606\"),\n  @(\"63\u00d752=3276\", \"30\u00d724=720\"),\n  @(\"33\u00d743=1419\", \"86\u00d768=5848\"),\n  @(\"40\u00d725=1000\", \"33\u00d764=2112\"),\n  @(\"85\u00d719=1615\", \"76\u00d741=3116\"),\n  @(\"63\u00d711=693\", \"45\u00d774=3330\"),\n  @(\"30\u00d790=2700\", \"88\u00d739=3432\"),\n  @(\"49\u00d733=1617\", \"56\u00d746=2576\"),\n  @(\"34\u00d761=2074\", \"84\u00d782=6888\"),\n  @(\"92\u00d765=5980\", \"42\u00d794=3948\"),\n  @(\"78\u00d761=4758\", \"62\u00d793=5766\"),\n  @(\"37\u00d741=1517\", \"46\u00d712=552\"),\n  @(\"56\u00d786=4816\", \"65\u00d772=4680\"),\n  @(\"38\u00d713=494\", \"71\u00d750=3550\"),\n  @(\"45\u00d766=2970\", \"39\u00d773=2847\"),\n  @(\"93\u00d723=2139\", \"81\u00d748=3888\")\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
